# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" sheet (detailed fund holdings) right before the
#    existing "总计" (summary) sheet.
# 2. Insert a new summary row at the top of "总计" for the 2022-Q1 quarter,
#    pushing the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: new "2022-Q1" sheet with per-fund holding detail
# ---------------------------------------------------------------------

$ws = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$ws.Name = "2022-Q1"

# NOTE: worksheet handles obtained *before* Worksheets.Add()/rename can end
# up pointing at the wrong sheet afterwards (this runtime resolves an
# already-fetched handle positionally, and Add() shifts everyone after the
# insertion point down by one slot) - so every sheet handle used below is
# (re)fetched by name only after all insert/rename calls that could move it.

# Style donor: any of the existing per-fund detail sheets carries the same
# header / row-number styling we want to replicate. (Deliberately NOT
# "总计" - copying from the sheet used as the new sheet's "Before" anchor
# trips a related addressing quirk where the pasted style silently fails
# to land.)
$styleSrc = $wb.Worksheets.Item("2021-Q4")

# Header row, styled like the other quarter sheets' header row (bold,
# centered, thin box border).
$styleSrc.Range("B1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Row-number column (A) styling, copied the same way.
$styleSrc.Range("A2").Copy()
$ws.Range("A2:A12").PasteSpecial(-4122)

$rows = @(
    @(0,  "010936", "交银施罗德均衡成长一年持有期混合A", "96.27", "87.48", "4.47", "4.3033", 4),
    @(1,  "519704", "交银先进制造混合",                   "75.33", "86.31", "4.42", "3.3296", 5),
    @(2,  "008099", "广发价值领先混合",                   "61.82", "83.88", "3.73", "2.3059", 9),
    @(3,  "005233", "广发睿毅领先混合",                   "40.39", "63.34", "4.50", "1.8176", 7),
    @(4,  "009402", "交银施罗德启明混合",                 "51.54", "82.62", "3.44", "1.7730", 6),
    @(5,  "001763", "广发多策略灵活配置混合",             "20.27", "69.36", "4.78", "0.9689", 8),
    @(6,  "270001", "广发聚富混合",                       "19.90", "73.54", "4.72", "0.9393", 8),
    @(7,  "519767", "交银施罗德科技创新灵活配置混合",     "4.18",  "90.38", "3.93", "0.1643", 6),
    @(8,  "010937", "交银施罗德均衡成长一年持有期混合C", "2.50",  "87.48", "4.47", "0.1118", 4),
    @(9,  "005104", "富荣福康混合A",                      "0.08",  "87.88", "3.08", "0.0025", 4),
    @(10, "005105", "富荣福康混合C",                      "0.04",  "87.88", "3.08", "0.0012", 4)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value  = $row[0]
    $ws.Cells.Item($r, 2).Value  = "'" + $row[1]
    $ws.Cells.Item($r, 3).Value  = $row[2]
    $ws.Cells.Item($r, 4).Value  = "'" + $row[3]
    $ws.Cells.Item($r, 5).Value  = "'" + $row[4]
    $ws.Cells.Item($r, 6).Value  = "'" + $row[5]
    $ws.Cells.Item($r, 7).Value  = "'" + $row[6]
    $ws.Cells.Item($r, 8).Value  = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Part 2: add the 2022-Q1 summary line to the "总计" sheet
# ---------------------------------------------------------------------

# Re-fetch by name now that the sheet collection has settled (see note
# above) rather than reusing a handle captured before Worksheets.Add().
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A2:D2").ClearFormats()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 11
$totalSheet.Range("D2").Value = 15.72

# Column A is a plain 0-based row index - renumber the rows that got
# pushed down (2021-Q4 .. 2020-Q4, now rows 3..7) to 1..5.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5

Write-Host "2022-Q1 data added"
